$p = $ppt.ActivePresentation
$s = $p.Slides.Item(27)

function FindByName($name) {
    for ($i = 1; $i -le $s.Shapes.Count; $i++) {
        $sh = $s.Shapes.Item($i)
        if ($sh.Name -eq $name) { return $sh }
    }
    return $null
}

$shp = FindByName("Picture 4")
Write-Output ("Picture 4 starts at pos=" + $shp.ZOrderPosition)
$shp.ZOrder(1)
$shp = FindByName("Picture 4")
Write-Output ("After ZOrder(1): pos=" + $shp.ZOrderPosition)
